$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ref, $value, $forceText) {
    $cell = $ws.Range($ref)
    if ($forceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}

$updates = @(
    @{ Ref = "D2"; Value = "67.103.82"; ForceText = $false },
    @{ Ref = "E2"; Value = "  +1.53%  "; ForceText = $false },
    @{ Ref = "D3"; Value = "3.896.30"; ForceText = $false },
    @{ Ref = "E3"; Value = "  +2.61%  "; ForceText = $false },
    @{ Ref = "D4"; Value = "0.999"; ForceText = $true },
    @{ Ref = "E4"; Value = "  +0.04%  "; ForceText = $false },
    @{ Ref = "D5"; Value = "466.71"; ForceText = $true },
    @{ Ref = "E5"; Value = "  +8.88%  "; ForceText = $false },
    @{ Ref = "D6"; Value = "144.77"; ForceText = $true },
    @{ Ref = "E6"; Value = "  +4.28%  "; ForceText = $false },
    @{ Ref = "D7"; Value = "0.627"; ForceText = $true },
    @{ Ref = "E7"; Value = "  +0.96%  "; ForceText = $false },
    @{ Ref = "E8"; Value = "  -0.04%  "; ForceText = $false },
    @{ Ref = "D9"; Value = "0.739"; ForceText = $true },
    @{ Ref = "E9"; Value = "  -0.24%  "; ForceText = $false },
    @{ Ref = "E10"; Value = "  +7.80%  "; ForceText = $false },
    @{ Ref = "D11"; Value = "0.0000338"; ForceText = $true },
    @{ Ref = "E11"; Value = "  +7.48%  "; ForceText = $false },
    @{ Ref = "D12"; Value = "43.06"; ForceText = $true },
    @{ Ref = "D13"; Value = "10.41"; ForceText = $true },
    @{ Ref = "E13"; Value = "  -1.31%  "; ForceText = $false },
    @{ Ref = "D14"; Value = "4.522.50"; ForceText = $false },
    @{ Ref = "E14"; Value = "  +3.07%  "; ForceText = $false },
    @{ Ref = "D15"; Value = "15.33"; ForceText = $true },
    @{ Ref = "E15"; Value = "  +1.72%  "; ForceText = $false },
    @{ Ref = "D16"; Value = "3.903.99"; ForceText = $false },
    @{ Ref = "E16"; Value = "  +2.81%  "; ForceText = $false },
    @{ Ref = "E17"; Value = "  -0.45%  "; ForceText = $false },
    @{ Ref = "D18"; Value = "19.95"; ForceText = $true },
    @{ Ref = "E18"; Value = "  -0.35%  "; ForceText = $false },
    @{ Ref = "E19"; Value = "  +3.02%  "; ForceText = $false },
    @{ Ref = "D20"; Value = "67.318.23"; ForceText = $false },
    @{ Ref = "E20"; Value = "  +1.72%  "; ForceText = $false },
    @{ Ref = "D21"; Value = "432.12"; ForceText = $true },
    @{ Ref = "E21"; Value = "  +6.24%  "; ForceText = $false },
    @{ Ref = "D22"; Value = "14.78"; ForceText = $true },
    @{ Ref = "E22"; Value = "  -2.79%  "; ForceText = $false },
    @{ Ref = "E23"; Value = "  +3.44%  "; ForceText = $false },
    @{ Ref = "D24"; Value = "88.58"; ForceText = $true },
    @{ Ref = "E24"; Value = "  +3.86%  "; ForceText = $false },
    @{ Ref = "D25"; Value = "38.86"; ForceText = $true },
    @{ Ref = "E25"; Value = "  +5.63%  "; ForceText = $false },
    @{ Ref = "D26"; Value = "3.53"; ForceText = $true },
    @{ Ref = "E26"; Value = "  +6.84%  "; ForceText = $false },
    @{ Ref = "D27"; Value = "5.74"; ForceText = $true },
    @{ Ref = "E27"; Value = "  +6.05%  "; ForceText = $false },
    @{ Ref = "D28"; Value = "10.11"; ForceText = $true },
    @{ Ref = "E28"; Value = "  +2.26%  "; ForceText = $false },
    @{ Ref = "D29"; Value = "9.59"; ForceText = $true },
    @{ Ref = "E29"; Value = "  -4.12%  "; ForceText = $false },
    @{ Ref = "D30"; Value = "736.77"; ForceText = $true },
    @{ Ref = "E30"; Value = "  +3.79%  "; ForceText = $false },
    @{ Ref = "D31"; Value = "13.73"; ForceText = $true },
    @{ Ref = "E31"; Value = "  -1.34%  "; ForceText = $false },
    @{ Ref = "E32"; Value = "  -1.09%  "; ForceText = $false },
    @{ Ref = "E33"; Value = "  +0.47%  "; ForceText = $false },
    @{ Ref = "D34"; Value = "43.20"; ForceText = $true },
    @{ Ref = "E34"; Value = "  +5.19%  "; ForceText = $false },
    @{ Ref = "D35"; Value = "0.157"; ForceText = $true },
    @{ Ref = "E35"; Value = "  +4.52%  "; ForceText = $false },
    @{ Ref = "D36"; Value = "58.14"; ForceText = $true },
    @{ Ref = "E36"; Value = "  +2.68%  "; ForceText = $false },
    @{ Ref = "E37"; Value = "  -0.20%  "; ForceText = $false },
    @{ Ref = "D38"; Value = "0.0₃0790"; ForceText = $false },
    @{ Ref = "E38"; Value = "  +15.40%  "; ForceText = $false },
    @{ Ref = "D39"; Value = "5.39"; ForceText = $true },
    @{ Ref = "E39"; Value = "  -5.81%  "; ForceText = $false },
    @{ Ref = "D40"; Value = "3.25"; ForceText = $true },
    @{ Ref = "E40"; Value = "  +13.56%  "; ForceText = $false },
    @{ Ref = "E41"; Value = "  -0.03%  "; ForceText = $false },
    @{ Ref = "E42"; Value = "  -0.65%  "; ForceText = $false },
    @{ Ref = "E43"; Value = "  +4.55%  "; ForceText = $false },
    @{ Ref = "E45"; Value = "  +5.52%  "; ForceText = $false },
    @{ Ref = "E46"; Value = "  +4.54%  "; ForceText = $false },
    @{ Ref = "E47"; Value = "  -6.63%  "; ForceText = $false },
    @{ Ref = "D48"; Value = "3.40"; ForceText = $true },
    @{ Ref = "E48"; Value = "  +0.27%  "; ForceText = $false },
    @{ Ref = "D49"; Value = "2.90"; ForceText = $true },
    @{ Ref = "E49"; Value = "  +3.26%  "; ForceText = $false },
    @{ Ref = "D50"; Value = "3.12"; ForceText = $true },
    @{ Ref = "E50"; Value = "  -0.52%  "; ForceText = $false },
    @{ Ref = "E51"; Value = "  +0.68%  "; ForceText = $false }
)

foreach ($u in $updates) {
    Set-CellText $u.Ref $u.Value $u.ForceText
}
